$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1022.55554
$ws.Range("I88").Value = 1016.6667
$ws.Range("J88").Value = 1023.73334
$ws.Range("K88").Value = 1016.6667
$ws.Range("L88").Value = 1023.73334
$ws.Range("M88").Value = -610.6667
$ws.Range("N88").Value = -1835.73334

$ws.Range("H91").Value = 1022.55554
$ws.Range("I91").Value = 1016.6667
$ws.Range("J91").Value = 1023.73334
$ws.Range("K91").Value = 1016.6667
$ws.Range("L91").Value = 1023.73334
$ws.Range("M91").Value = 387.3333
$ws.Range("N91").Value = -3831.73334

$ws.Range("H100").Value = 19610190
$ws.Range("I100").Value = 37039160
$ws.Range("J100").Value = 2598.75
$ws.Range("K100").Value = 37039160
$ws.Range("L100").Value = 2598.75
$ws.Range("M100").Value = -37038619
$ws.Range("N100").Value = -3680.75

$ws.Range("H131").Value = 4696.38
$ws.Range("I131").Value = 1229.75
$ws.Range("J131").Value = 4997.826
$ws.Range("K131").Value = 3689.25
$ws.Range("L131").Value = 14993.478
$ws.Range("M131").Value = 1350.75
$ws.Range("N131").Value = -25073.478

$ws.Range("H138").Value = 2953.8872
$ws.Range("I138").Value = 2133.8572
$ws.Range("J138").Value = 3193.0625
$ws.Range("K138").Value = 6401.571599999999
$ws.Range("L138").Value = 9579.1875
$ws.Range("M138").Value = -1261.571599999999
$ws.Range("N138").Value = -19859.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3165.6667
$ws.Range("J88").Value = 3165.6667
$ws.Range("L88").Value = 3165.6667
$ws.Range("N88").Value = -3977.6667

$ws.Range("H91").Value = 3165.6667
$ws.Range("J91").Value = 3165.6667
$ws.Range("L91").Value = 3165.6667
$ws.Range("N91").Value = -5973.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H134").Value = 2389.4333
$ws.Range("I134").Value = 1234.3182
$ws.Range("J134").Value = 5566
$ws.Range("K134").Value = 3702.9546
$ws.Range("L134").Value = 16698
$ws.Range("M134").Value = -1167.9546
$ws.Range("N134").Value = -21768

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2548.32
$ws.Range("I31").Value = 2166.3333
$ws.Range("J31").Value = 3530.5715
$ws.Range("K31").Value = 2166.3333
$ws.Range("L31").Value = 3530.5715
$ws.Range("M31").Value = -1871.3333
$ws.Range("N31").Value = -4120.5715

$ws.Range("H34").Value = 2548.32
$ws.Range("I34").Value = 2166.3333
$ws.Range("J34").Value = 3530.5715
$ws.Range("K34").Value = 2166.3333
$ws.Range("L34").Value = 3530.5715
$ws.Range("M34").Value = -1964.3333
$ws.Range("N34").Value = -3934.5715

$ws.Range("H86").Value = 3254.7693
$ws.Range("J86").Value = 2962.6
$ws.Range("L86").Value = 2962.6
$ws.Range("N86").Value = -5208.6

$ws.Range("H89").Value = 3254.7693
$ws.Range("J89").Value = 2962.6
$ws.Range("L89").Value = 14813
$ws.Range("N89").Value = -26045

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1613.9375
$ws.Range("I5").Value = 742.8
$ws.Range("J5").Value = 3065.8333
$ws.Range("K5").Value = 2228.4
$ws.Range("L5").Value = 9197.499899999999
$ws.Range("M5").Value = -2116.4
$ws.Range("N5").Value = -9421.499899999999

$ws.Range("H92").Value = 1044.2
$ws.Range("I92").Value = 811
$ws.Range("J92").Value = 1199.6666
$ws.Range("K92").Value = 2433
$ws.Range("L92").Value = 3598.9998
$ws.Range("M92").Value = -1185
$ws.Range("N92").Value = -6094.9998

$ws.Range("H122").Value = 631.7
$ws.Range("I122").Value = 270.5
$ws.Range("J122").Value = 1173.5
$ws.Range("K122").Value = 2434.5
$ws.Range("L122").Value = 10561.5
$ws.Range("M122").Value = 15.5
$ws.Range("N122").Value = -15461.5

$ws.Range("H131").Value = 2902.8103
$ws.Range("J131").Value = 2948.1228
$ws.Range("L131").Value = 8844.368399999999
$ws.Range("N131").Value = -18924.3684

$ws.Range("H135").Value = 1613.9375
$ws.Range("I135").Value = 742.8
$ws.Range("J135").Value = 3065.8333
$ws.Range("K135").Value = 6685.2
$ws.Range("L135").Value = 27592.4997
$ws.Range("M135").Value = -4150.2
$ws.Range("N135").Value = -32662.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2299.389
$ws.Range("I132").Value = 1493.2307
$ws.Range("J132").Value = 4395.4
$ws.Range("K132").Value = 4479.6921
$ws.Range("L132").Value = 13186.2
$ws.Range("M132").Value = -1949.6921
$ws.Range("N132").Value = -18246.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 144929
$ws.Range("I7").Value = 168083.17
$ws.Range("J7").Value = 6004
$ws.Range("K7").Value = 168083.17
$ws.Range("L7").Value = 6004
$ws.Range("M7").Value = -167971.17
$ws.Range("N7").Value = -6228

$ws.Range("H82").Value = 1650
$ws.Range("I82").Value = 1517
$ws.Range("J82").Value = 1699.875
$ws.Range("K82").Value = 1517
$ws.Range("L82").Value = 1699.875
$ws.Range("M82").Value = -1156
$ws.Range("N82").Value = -2421.875

$ws.Range("H85").Value = 1650
$ws.Range("I85").Value = 1517
$ws.Range("J85").Value = 1699.875
$ws.Range("K85").Value = 1517
$ws.Range("L85").Value = 1699.875
$ws.Range("M85").Value = -269
$ws.Range("N85").Value = -4195.875

$ws.Range("H126").Value = 144929
$ws.Range("I126").Value = 168083.17
$ws.Range("J126").Value = 6004
$ws.Range("K126").Value = 504249.51
$ws.Range("L126").Value = 18012
$ws.Range("M126").Value = -501779.51
$ws.Range("N126").Value = -22952

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1910.0625
$ws.Range("I81").Value = 2200.125
$ws.Range("J81").Value = 1620
$ws.Range("K81").Value = 4400.25
$ws.Range("L81").Value = 3240
$ws.Range("M81").Value = -3339.25
$ws.Range("N81").Value = -5362

$ws.Range("H84").Value = 1910.0625
$ws.Range("I84").Value = 2200.125
$ws.Range("J84").Value = 1620
$ws.Range("K84").Value = 22001.25
$ws.Range("L84").Value = 16200
$ws.Range("M84").Value = -16697.25
$ws.Range("N84").Value = -26808

$ws.Range("H86").Value = 25162.5
$ws.Range("J86").Value = 25162.5
$ws.Range("L86").Value = 25162.5
$ws.Range("N86").Value = -27408.5

$ws.Range("H89").Value = 25162.5
$ws.Range("J89").Value = 25162.5
$ws.Range("L89").Value = 125812.5
$ws.Range("N89").Value = -137044.5

$ws.Range("H132").Value = 1749.1515
$ws.Range("I132").Value = 1441.4445
$ws.Range("K132").Value = 4324.333500000001
$ws.Range("M132").Value = -1794.333500000001

$ws.Range("H136").Value = 6668779
$ws.Range("I136").Value = 9009549
$ws.Range("J136").Value = 6585.385
$ws.Range("K136").Value = 27028647
$ws.Range("L136").Value = 19756.155
$ws.Range("M136").Value = -27026097
$ws.Range("N136").Value = -24856.155
